$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "alamat" column: shift no_telepon/email left into C1/D1,
# and clear the now-unused E1 header cell (keeping column widths/cols
# definitions untouched, unlike a structural column delete).
$ws.Range("C1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("E1").Value()
$ws.Range("E1").ClearContents() | Out-Null

# Remove the long "Alamat lengkap..." note that lived in row 2
$ws.Range("C2").ClearContents() | Out-Null

# Move the active selection to C6, matching the saved workbook state
$ws.Range("C6").Select() | Out-Null
